$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.735.50"
$ws.Range("E2").Value = "  +1.21%  "

$ws.Range("D3").Value = "3.739.22"
$ws.Range("E3").Value = "  -1.68%  "

$ws.Range("E4").Value = "  +0.05%  "

$ws.Range("D5").Value = "'601.62"
$ws.Range("E5").Value = "  +0.01%  "

$ws.Range("D6").Value = "'168.71"
$ws.Range("E6").Value = "  -2.09%  "

$ws.Range("D7").Value = "3.736.68"
$ws.Range("E7").Value = "  -1.68%  "

$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("E9").Value = "  +0.66%  "

$ws.Range("E10").Value = "  +2.46%  "

$ws.Range("D11").Value = "'6.35"
$ws.Range("E11").Value = "  +2.32%  "

$ws.Range("E12").Value = "  -1.15%  "

$ws.Range("D13").Value = "'38.18"
$ws.Range("E13").Value = "  -1.68%  "

$ws.Range("D14").Value = "'0.0000245"
$ws.Range("E14").Value = "  +0.24%  "

$ws.Range("D15").Value = "4.364.44"
$ws.Range("E15").Value = "  -1.66%  "

$ws.Range("D16").Value = "3.738.48"
$ws.Range("E16").Value = "  -1.63%  "

$ws.Range("D17").Value = "68.746.66"
$ws.Range("E17").Value = "  +1.26%  "

$ws.Range("D18").Value = "'7.31"
$ws.Range("E18").Value = "  +0.44%  "

$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("D20").Value = "'17.14"
$ws.Range("E20").Value = "  -0.76%  "

$ws.Range("D21").Value = "'10.70"
$ws.Range("E21").Value = "  +16.16%  "

$ws.Range("D22").Value = "'494.05"
$ws.Range("E22").Value = "  +0.01%  "

$ws.Range("E23").Value = "  -2.00%  "

$ws.Range("D24").Value = "'85.29"
$ws.Range("E24").Value = "  -0.58%  "

$ws.Range("E25").Value = "  -0.60%  "

$ws.Range("E26").Value = "  -3.62%  "

$ws.Range("D27").Value = "'12.41"
$ws.Range("E27").Value = "  +0.18%  "

$ws.Range("D28").Value = "'10.14"
$ws.Range("E28").Value = "  -0.74%  "

$ws.Range("E29").Value = "  -0.08%  "

$ws.Range("D30").Value = "'2.56"
$ws.Range("E30").Value = "  +4.64%  "

$ws.Range("D31").Value = "'2.97"
$ws.Range("E31").Value = "  +0.06%  "

$ws.Range("D32").Value = "'7.91"
$ws.Range("E32").Value = "  +0.53%  "

$ws.Range("D33").Value = "'31.69"
$ws.Range("E33").Value = "  -3.99%  "

$ws.Range("D34").Value = "3.883.85"
$ws.Range("E34").Value = "  -1.56%  "

$ws.Range("D35").Value = "3.673.15"
$ws.Range("E35").Value = "  -1.75%  "

$ws.Range("E36").Value = "  -1.71%  "

$ws.Range("D37").Value = "'1.00"
$ws.Range("E37").Value = "  +0.09%  "

$ws.Range("E38").Value = "  -0.49%  "

$ws.Range("D39").Value = "'5.82"
$ws.Range("E39").Value = "  -0.37%  "

$ws.Range("E40").Value = "  +0.58%  "

$ws.Range("E41").Value = "  -1.94%  "

$ws.Range("D42").Value = "'436.70"
$ws.Range("E42").Value = "  -5.82%  "

$ws.Range("D43").Value = "'48.92"
$ws.Range("E43").Value = "  -0.24%  "

$ws.Range("E44").Value = "  -1.68%  "

$ws.Range("D45").Value = "'2.86"
$ws.Range("E45").Value = "  +0.15%  "

$ws.Range("E46").Value = "  +0.73%  "

$ws.Range("E47").Value = "  +0.02%  "

$ws.Range("D48").Value = "'40.64"
$ws.Range("E48").Value = "  -1.09%  "

$ws.Range("D49").Value = "'141.17"
$ws.Range("E49").Value = "  +1.02%  "

$ws.Range("D50").Value = "2.794.41"
$ws.Range("E50").Value = "  -1.80%  "
